$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Optieorder" (Kop3 heading, followed by a separate run " aanpassen")
#    becomes "Wijze van invoer van optietransactie" (heading keeps its
#    trailing " aanpassen").
# ------------------------------------------------------------------
$headingIdx = -1
$priceIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Optieorder aanpassen") {
        $headingIdx = $i
    }
    if ($t -eq "Optieprijsberekening en optieprijs aanpassen") {
        $priceIdx = $i
    }
}

$heading = $d.Paragraphs.Item($headingIdx)
$headingStart = $heading.Range.Start
$replaceRange = $d.Range($headingStart, $headingStart + ("Optieorder").Length)
$replaceRange.Find.Execute("Optieorder", $false, $true, $false, $false, $false, $true, 0, $false, "Wijze van invoer van optietransactie", 1)

# ------------------------------------------------------------------
# 2) After the "Optieprijsberekening en optieprijs aanpassen" heading
#    (and the blank paragraph that already follows it) insert a new
#    remark paragraph plus a trailing blank paragraph.
# ------------------------------------------------------------------
$priceHeading = $d.Paragraphs.Item($priceIdx)
$blankAfterPrice = $priceHeading.Next()

$newPara = $blankAfterPrice.Range.InsertParagraphAfter()
$insertedIdx = $priceIdx + 2
$d.Paragraphs.Item($insertedIdx).Range.Text = "Moet trouwens zijn " + [char]0x201C + "Optietransactie uitvoeren" + [char]0x201D + ". Er wordt voor opties geen order ingelegd in het pakket, transactie wordt direct uitgevoerd."
$d.Paragraphs.Item($insertedIdx).Range.InsertParagraphAfter()
